# Applies the "Updated cryptos list" data refresh (rates + row-26/27 swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.597.53'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '2.067.64'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'242.07"
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").Value = "'0.673"
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'52.83"
$ws.Range("E8").Value = '  -6.62%  '
$ws.Range("D9").Value = "'58.84"
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("D10").Value = "'0.360"
$ws.Range("E10").Value = '  -6.06%  '
$ws.Range("D11").Value = "'0.0751"
$ws.Range("E11").Value = '  -3.50%  '
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = "'0.895"
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("E14").Value = '  -9.10%  '
$ws.Range("D15").Value = '2.369.90'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = '  -5.00%  '
$ws.Range("D17").Value = '2.015.16'
$ws.Range("E17").Value = '  -1.19%  '
$ws.Range("D18").Value = '36.542.94'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("D19").Value = "'16.37"
$ws.Range("E19").Value = '  -13.83%  '
$ws.Range("D20").Value = "'71.74"
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("D21").Value = '0.0₃0865'
$ws.Range("E21").Value = '  -2.71%  '
$ws.Range("D22").Value = "'5.28"
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").Value = "'235.84"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = '  -4.71%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = "'2.13"
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = "'9.27"
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("D28").Value = "'163.22"
$ws.Range("E28").Value = '  -4.82%  '
$ws.Range("D29").Value = "'20.52"
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").Value = "'4.58"
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").Value = "'0.0597"
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("E38").Value = '  -6.40%  '
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("E40").Value = '  -4.92%  '
$ws.Range("D41").Value = "'4.86"
$ws.Range("E41").Value = '  -6.01%  '
$ws.Range("D42").Value = "'0.0216"
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("D43").Value = "'1.13"
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("E44").Value = '  -6.18%  '
$ws.Range("D45").Value = "'94.02"
$ws.Range("E45").Value = '  -3.69%  '
$ws.Range("D46").Value = '1.395.34'
$ws.Range("E46").Value = '  +8.80%  '
$ws.Range("D47").Value = "'15.59"
$ws.Range("E47").Value = '  -9.34%  '
$ws.Range("D48").Value = "'7.33"
$ws.Range("E48").Value = '  +7.99%  '
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("D50").Value = "'2.86"
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("D51").Value = '2.255.28'
$ws.Range("E51").Value = '  +0.98%  '
